$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45936

$ws.Range("B2").Value = 45.08
$ws.Range("C2").Value = 42.42
$ws.Range("D2").Value = 43.05
$ws.Range("E2").Value = 45.23
$ws.Range("F2").Value = 46.81
$ws.Range("G2").Value = 55.81
$ws.Range("H2").Value = 82.06
$ws.Range("I2").Value = 112.32
$ws.Range("J2").Value = 114.61
$ws.Range("K2").Value = 106.67
$ws.Range("L2").Value = 70.77
$ws.Range("M2").Value = 46.97
$ws.Range("N2").Value = 44.44
$ws.Range("O2").Value = 41.44
$ws.Range("P2").Value = 36.31
$ws.Range("Q2").Value = 40.33
$ws.Range("R2").Value = 43.78
$ws.Range("S2").Value = 68.62
$ws.Range("T2").Value = 98.93000000000001
$ws.Range("U2").Value = 137.64
$ws.Range("V2").Value = 163.34
$ws.Range("W2").Value = 169.24
$ws.Range("X2").Value = 125.74
$ws.Range("Y2").Value = 103.13
$ws.Range("Z2").Value = 78.53

$ws.Range("AB2").Value = 140.36
$ws.Range("AD2").Value = 166.29
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 118.28
$ws.Range("AG2").Value = "0h-17h"
